# Applies the OOXML diff:
#  - Column B and E width: 40 -> 30
#  - F3: "10:36:00" -> "8:54:00"
#  - F4: "11:40:13" -> "12:08:00"
#  - D5: "-" -> "9LVB3XPUL11"
#  - E5: "Version 15.0(TTC_20140605)FLO_DSGS7" -> "Version 15.2(4.0.55)E"
#  - F5: "11:35:00" -> "12:04:00"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Narrow columns B and E from 40 to 30 (stored XML "width" attribute).
# Excel's ColumnWidth property is in character units with a ~0.83 offset
# from the raw stored width at this workbook's default font, so use the
# character-unit value that round-trips to a stored width of exactly 30.
$ws.Columns.Item(2).ColumnWidth = 29.17
$ws.Columns.Item(5).ColumnWidth = 29.17

# Update uptime values (column F) - keep as text
$ws.Range("F3").Value = "8:54:00"
$ws.Range("F4").Value = "12:08:00"
$ws.Range("F5").Value = "12:04:00"

# Update serial_number and os_version for SW1 (row 5)
$ws.Range("D5").Value = "9LVB3XPUL11"
$ws.Range("E5").Value = "Version 15.2(4.0.55)E"
